# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Row (by sheet row number) -> new F-column value, for sheet "展览"
$updatesExhibition = @{
    2  = 1065
    3  = 348
    4  = 1460
    5  = 8668
    7  = 487
    8  = 637
    9  = 275
    11 = 5
    12 = 3516
    14 = 359
    16 = 1106
    17 = 144
    19 = 305
    20 = 192
    21 = 2238
    22 = 46
}

# Row (by sheet row number) -> new F-column value, for sheet "全部类型"
$updatesAll = @{
    2  = 1065
    3  = 348
    4  = 1460
    5  = 8668
    7  = 487
    8  = 637
    9  = 275
    11 = 5
    12 = 3516
    14 = 359
    16 = 1106
    17 = 144
    19 = 305
    20 = 192
    21 = 2238
    23 = 46
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $updatesExhibition[$row]
}

$wsAll = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAll.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $updatesAll[$row]
}
